$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Self-assessment checklist")
$ws.Rows.Item(40).Delete()
